$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

$ws.Cells.Item($row, 1).Value = "CompaNanny"
$ws.Cells.Item($row, 2).Value = "CompaNanny Statenkwartier KDV"
$ws.Cells.Item($row, 3).Value = "KDV"

# Column D holds a plain text date string ("2023-05-04"), not a real date.
# Force text formatting first so Excel doesn't auto-convert it to a date
# serial number, then strip the style back off so the cell matches the
# unstyled cells around it.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2023-05-04"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
